# Update "想去人数" (F column) values across the sheets as per the
# regenerated gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 43
$ws1.Range("F3").Value = 26421
$ws1.Range("F4").Value = 581
$ws1.Range("F6").Value = 596
$ws1.Range("F7").Value = 175
$ws1.Range("F8").Value = 537
$ws1.Range("F9").Value = 232
$ws1.Range("F10").Value = 353
$ws1.Range("F11").Value = 223
$ws1.Range("F12").Value = 187
$ws1.Range("F15").Value = 50
$ws1.Range("F16").Value = 382
$ws1.Range("F18").Value = 1514
$ws1.Range("F19").Value = 188
$ws1.Range("F20").Value = 33

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 228
$ws2.Range("F6").Value = 186
$ws2.Range("F8").Value = 110
$ws2.Range("F9").Value = 110
$ws2.Range("F10").Value = 434

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5013
$ws3.Range("F3").Value = 214

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 43
$ws4.Range("F3").Value = 5013
$ws4.Range("F4").Value = 214
$ws4.Range("F5").Value = 26421
$ws4.Range("F6").Value = 581
$ws4.Range("F9").Value = 228
$ws4.Range("F10").Value = 596
$ws4.Range("F13").Value = 175
$ws4.Range("F14").Value = 186
$ws4.Range("F15").Value = 186
$ws4.Range("F17").Value = 110
$ws4.Range("F18").Value = 110
$ws4.Range("F19").Value = 434
$ws4.Range("F20").Value = 537
$ws4.Range("F22").Value = 232
$ws4.Range("F23").Value = 353
$ws4.Range("F24").Value = 223
$ws4.Range("F25").Value = 187
$ws4.Range("F29").Value = 50
$ws4.Range("F32").Value = 382
$ws4.Range("F35").Value = 1514
$ws4.Range("F36").Value = 188
$ws4.Range("F38").Value = 33
